$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet stores runs/balls/fours/sixes as text-formatted numbers
# (numberStoredAsText). Go through a text formula + paste-values so the
# resulting cells stay text (shared-string) rather than becoming real
# numbers, matching the rest of the sheet.

# Row 2 (runs, balls, fours, sixes)
$ws.Range("C2").Formula = "=""11"""
$ws.Range("C2").Copy() | Out-Null
$ws.Range("C2").PasteSpecial(-4163) | Out-Null

$ws.Range("D2").Formula = "=""13"""
$ws.Range("D2").Copy() | Out-Null
$ws.Range("D2").PasteSpecial(-4163) | Out-Null

$ws.Range("E2").Formula = "=""1"""
$ws.Range("E2").Copy() | Out-Null
$ws.Range("E2").PasteSpecial(-4163) | Out-Null

$ws.Range("F2").Formula = "=""0"""
$ws.Range("F2").Copy() | Out-Null
$ws.Range("F2").PasteSpecial(-4163) | Out-Null

# Row 3
$ws.Range("C3").Formula = "=""47"""
$ws.Range("C3").Copy() | Out-Null
$ws.Range("C3").PasteSpecial(-4163) | Out-Null

$ws.Range("D3").Formula = "=""39"""
$ws.Range("D3").Copy() | Out-Null
$ws.Range("D3").PasteSpecial(-4163) | Out-Null

$ws.Range("E3").Formula = "=""1"""
$ws.Range("E3").Copy() | Out-Null
$ws.Range("E3").PasteSpecial(-4163) | Out-Null

$ws.Range("F3").Formula = "=""3"""
$ws.Range("F3").Copy() | Out-Null
$ws.Range("F3").PasteSpecial(-4163) | Out-Null

# Row 4
$ws.Range("C4").Formula = "=""1"""
$ws.Range("C4").Copy() | Out-Null
$ws.Range("C4").PasteSpecial(-4163) | Out-Null

$ws.Range("D4").Formula = "=""2"""
$ws.Range("D4").Copy() | Out-Null
$ws.Range("D4").PasteSpecial(-4163) | Out-Null

$ws.Range("E4").Formula = "=""0"""
$ws.Range("E4").Copy() | Out-Null
$ws.Range("E4").PasteSpecial(-4163) | Out-Null

$ws.Range("F4").Formula = "=""0"""
$ws.Range("F4").Copy() | Out-Null
$ws.Range("F4").PasteSpecial(-4163) | Out-Null

$excel.CutCopyMode = $false
